# Insert a new data row above the current row 221 (pushes the existing
# rows 221-333 down to 222-334, preserving their formatting/values), then
# populate the newly-inserted row 221 with its own record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("221:221").Insert()

$ws.Range("A221").Value = 4
$ws.Range("B221").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C221").Value = "Los Lagos"
$ws.Range("D221").Value = 44813
$ws.Range("E221").Value = 10
$ws.Range("F221").Value = 100112037
$ws.Range("G221").Value = "Cebollín"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 180
$ws.Range("K221").Value = 9000
$ws.Range("L221").Value = 9000
$ws.Range("M221").Value = 9000
$ws.Range("N221").Value = "$/paquete 36 unidades"
$ws.Range("O221").Value = "Región Metropolitana"
$ws.Range("P221").Value = 250
$ws.Range("Q221").Value = 36
$ws.Range("R221").Value = "Hortaliza"
